$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. SupIm sheet: add 11 new timeseries rows (rows 4-14), mirroring row 3's
#    values/pattern (hour index in column A incrementing, same B/C/D values).
# ---------------------------------------------------------------------------
$supim = $wb.Worksheets.Item("SupIm")

for ($r = 4; $r -le 14; $r++) {
    $supim.Range("A$r").Value = $r - 2
    $supim.Range("B$r").Value = 0.481
    $supim.Range("C$r").Value = 0.3
    $supim.Range("D$r").Value = 0.207
}

# Copy formatting from the template row (row 3) down across the new rows so
# number formats / fills match exactly, same as a fill-down / paste-format.
$supim.Range("A3:D3").Copy()
$supim.Range("A4:D14").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Merge the Process sheet's conditional formatting rules that covered
#    "A12:C13 A11 C11" and "B11" separately into one rule over "A11:C13".
# ---------------------------------------------------------------------------
$process = $wb.Worksheets.Item("Process")

$bigRuleRange = $process.Range("A12:C13,A11,C11")
$smallRuleRange = $process.Range("B11")

$bigRule = $bigRuleRange.FormatConditions.Item(1)
$smallRule = $smallRuleRange.FormatConditions.Item(1)

$smallRule.Delete()
$bigRule.ModifyAppliesToRange($process.Range("A11:C13"))
$bigRule.SetFirstPriority()

# ---------------------------------------------------------------------------
# 3. Update the active sheet / selection state: SupIm becomes the active
#    (selected) sheet with J15 selected, Process loses tabSelected.
# ---------------------------------------------------------------------------
$supim.Activate()
$supim.Range("J15").Select()
